# ZBP_05_kontakt_s_lidmi.xlsx weekly refresh
#
# - adds two new weekly columns ("15.–21. 2. 2021" and "22.–28. 2. 2021") to both
#   the "data" sheet (as AV/AW, continuing after existing AU) and the "pocetR"
#   sheet (as AU/AV, continuing after existing AT)
# - bumps the "aktualizace" date in the two trailing footnote cells from
#   23. 2. 2021 to 9. 3. 2021
# - applies a handful of small retrospective sample-size corrections in
#   "pocetR" rows 20-21 that shipped with the same update

$wb = $excel.ActiveWorkbook

$wsData   = $wb.Worksheets.Item("data")
$wsPocetR = $wb.Worksheets.Item("pocetR")

$label1 = "15.–21. 2. 2021"
$label2 = "22.–28. 2. 2021"

# ---------------------------------------------------------------------------
# Sheet "data": new columns AV (48) and AW (49), continuing after AU (47)
# ---------------------------------------------------------------------------

# Header row - copy formatting from the last existing header cell, then set values
$wsData.Range("AU1").Copy()
$wsData.Range("AV1:AW1").PasteSpecial(-4122)  # xlPasteFormats
$wsData.Range("AV1").Value = $label1
$wsData.Range("AW1").Value = $label2

# Data rows 2-61: AV/AW numeric values
$sheet1NewCols = @{
    2 = @(12, 12)
    3 = @(6, 6)
    4 = @(5, 5)
    5 = @(14, 14)
    6 = @(17, 16.5)
    7 = @(7, 7.5)
    8 = @(15, 14.5)
    9 = @(11.5, 12)
    10 = @(11.5, 11.5)
    11 = @(16, 16)
    12 = @(11, 10.5)
    13 = @(9, 9.5)
    14 = @(13.5, 13.5)
    15 = @(10.5, 10.5)
    16 = @(12.5, 12.5)
    17 = @(11, 11.5)
    18 = @(12, 11.5)
    19 = @(12, 11.5)
    20 = @(21.5, 21.5)
    21 = @(13, 13.5)
    22 = @(8, 7)
    23 = @(10.5, 8)
    24 = @(8, 8)
    25 = @(10, 8)
    26 = @(4, 4)
    27 = @(8, 8)
    28 = @(5, 6)
    29 = @(5, 5)
    30 = @(10, 10)
    31 = @(5, 5)
    32 = @(5, 5)
    33 = @(7, 7)
    34 = @(5, 5)
    35 = @(6, 6)
    36 = @(6, 6)
    37 = @(6, 6)
    38 = @(6, 5)
    39 = @(15, 15)
    40 = @(10, 10)
    41 = @(5, 5)
    42 = @(4, 3)
    43 = @(7, 7)
    44 = @(6, 6)
    45 = @(3, 3)
    46 = @(6, 6)
    47 = @(5, 5)
    48 = @(4.5, 4.5)
    49 = @(6.5, 6.5)
    50 = @(4.5, 4.5)
    51 = @(4, 4)
    52 = @(5, 5)
    53 = @(4.5, 4.5)
    54 = @(5.5, 5)
    55 = @(4.5, 4.5)
    56 = @(5, 5)
    57 = @(4.5, 4.5)
    58 = @(7, 7)
    59 = @(4.5, 4.5)
    60 = @(4.5, 5)
    61 = @(4.5, 4)
}

foreach ($r in $sheet1NewCols.Keys) {
    $vals = $sheet1NewCols[$r]
    $wsData.Cells.Item($r, 48).Value = $vals[0]
    $wsData.Cells.Item($r, 49).Value = $vals[1]
}

# Footnote row 62: bump the "aktualizace" date
$wsData.Range("A62").Value = "Život během pandemie, Kontakt s lidmi, průměr celkově a ve skupinách, aktualizace 9. 3. 2021"

# ---------------------------------------------------------------------------
# Sheet "pocetR": new columns AU (47) and AV (48), continuing after AT (46)
# ---------------------------------------------------------------------------

$wsPocetR.Range("AT1").Copy()
$wsPocetR.Range("AU1:AV1").PasteSpecial(-4122)  # xlPasteFormats
$wsPocetR.Range("AU1").Value = $label1
$wsPocetR.Range("AV1").Value = $label2

$sheet2NewCols = @{
    2 = @(1915, 1924)
    3 = @(444, 449)
    4 = @(669, 670)
    5 = @(802, 805)
    6 = @(335, 341)
    7 = @(597, 598)
    8 = @(983, 985)
    9 = @(640, 645)
    10 = @(680, 682)
    11 = @(595, 597)
    12 = @(940, 945)
    13 = @(975, 979)
    14 = @(991, 1001)
    15 = @(443, 443)
    16 = @(232, 231)
    17 = @(249, 249)
    18 = @(669, 674)
    19 = @(91, 91)
    20 = @(146, 147)
    21 = @(120, 120)
}

foreach ($r in $sheet2NewCols.Keys) {
    $vals = $sheet2NewCols[$r]
    $wsPocetR.Cells.Item($r, 47).Value = $vals[0]
    $wsPocetR.Cells.Item($r, 48).Value = $vals[1]
}

# Row 22 ("Celkem" grand-total row) only gets the two new, still-empty marker
# cells, matching the blank cells used across the rest of that row.
$wsPocetR.Range("AT22").Copy()
$wsPocetR.Range("AU22:AV22").PasteSpecial(-4122)  # xlPasteFormats

# Small retrospective sample-size corrections shipped in the same update
$sheet2Corrections = @{
    "I20"  = 270
    "M20"  = 165
    "Q20"  = 92
    "AA20" = 107
    "AG20" = 190
    "AK20" = 139
    "AM20" = 135
    "AQ20" = 145
    "AS20" = 156
    "E21"  = 299
    "I21"  = 251
    "M21"  = 153
    "Q21"  = 86
    "AA21" = 89
    "AG21" = 146
    "AK21" = 96
    "AM21" = 95
    "AQ21" = 120
}

foreach ($addr in $sheet2Corrections.Keys) {
    $wsPocetR.Range($addr).Value = $sheet2Corrections[$addr]
}

# Footnote row 22: bump the "aktualizace" date
$wsPocetR.Range("A22").Value = "Život během pandemie, Kontakt s lidmi, velikost dotázaného souboru celkově a ve skupinách, aktualizace 9. 3. 2021"

Write-Host "edit complete"
